# Apply the robot-update numeric corrections to the VD (X) and CH (AB) canton
# columns across the Cases, Fatalities, Hospitalized and ICU sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("X32").Value = 135
$ws.Range("AB32").Value = 886
$ws.Range("X33").Value = 144
$ws.Range("AB33").Value = 1024
$ws.Range("X34").Value = 169
$ws.Range("AB34").Value = 1169
$ws.Range("X35").Value = 183
$ws.Range("AB35").Value = 1280
$ws.Range("X36").Value = 186
$ws.Range("AB36").Value = 1323
$ws.Range("X37").Value = 192
$ws.Range("AB37").Value = 1375
$ws.Range("X38").Value = 206
$ws.Range("AB38").Value = 1454
$ws.Range("X39").Value = 222
$ws.Range("AB39").Value = 1540
$ws.Range("X40").Value = 228
$ws.Range("AB40").Value = 1629
$ws.Range("AB41").Value = 1663

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("X39").Value = 2
$ws.Range("X40").Value = 2

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("X3").Value = -1
$ws.Range("AB3").Value = -3
$ws.Range("X4").Value = 2
$ws.Range("AB4").Value = -3
$ws.Range("X5").Value = 2
$ws.Range("AB5").Value = -25
$ws.Range("X6").Value = 2
$ws.Range("AB6").Value = -26
$ws.Range("X7").Value = 3
$ws.Range("AB7").Value = -27
$ws.Range("X8").Value = -1
$ws.Range("AB8").Value = -37
$ws.Range("X9").Value = -1
$ws.Range("AB9").Value = -40
$ws.Range("X10").Value = 0
$ws.Range("AB10").Value = -54
$ws.Range("X14").Value = 0
$ws.Range("AB14").Value = -58
$ws.Range("X17").Value = -2
$ws.Range("AB17").Value = -75
$ws.Range("X18").Value = -5
$ws.Range("AB18").Value = -80
$ws.Range("X19").Value = -7
$ws.Range("AB19").Value = -77
$ws.Range("X20").Value = -8
$ws.Range("AB20").Value = -79
$ws.Range("X21").Value = -6
$ws.Range("AB21").Value = -81
$ws.Range("X22").Value = -6
$ws.Range("AB22").Value = -82
$ws.Range("X23").Value = -5
$ws.Range("AB23").Value = -80
$ws.Range("X24").Value = -9
$ws.Range("AB24").Value = -78
$ws.Range("X25").Value = -11
$ws.Range("AB25").Value = -83
$ws.Range("X26").Value = -12
$ws.Range("AB26").Value = -83
$ws.Range("X27").Value = -13
$ws.Range("AB27").Value = -87
$ws.Range("X28").Value = -12
$ws.Range("AB28").Value = -85
$ws.Range("X29").Value = -12
$ws.Range("AB29").Value = -86
$ws.Range("X30").Value = -12
$ws.Range("AB30").Value = -87
$ws.Range("X31").Value = -11
$ws.Range("AB31").Value = -79
$ws.Range("X32").Value = -12
$ws.Range("AB32").Value = -78
$ws.Range("X33").Value = -13
$ws.Range("AB33").Value = -76
$ws.Range("X34").Value = -12
$ws.Range("AB34").Value = -83
$ws.Range("X35").Value = -13
$ws.Range("AB35").Value = -85
$ws.Range("X36").Value = -14
$ws.Range("AB36").Value = -85
$ws.Range("X37").Value = -12
$ws.Range("AB37").Value = -80
$ws.Range("X38").Value = -18
$ws.Range("AB38").Value = -83
$ws.Range("X39").Value = -18
$ws.Range("AB39").Value = -81
$ws.Range("X40").Value = -16
$ws.Range("AB40").Value = -86
$ws.Range("AB41").Value = -83

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("X7").Value = -1
$ws.Range("AB7").Value = -3
$ws.Range("X8").Value = -1
$ws.Range("AB8").Value = -4
$ws.Range("X9").Value = -1
$ws.Range("AB9").Value = -3
$ws.Range("X10").Value = -2
$ws.Range("AB10").Value = -8
$ws.Range("X11").Value = -2
$ws.Range("AB11").Value = -10
$ws.Range("X12").Value = -3
$ws.Range("AB12").Value = -13
$ws.Range("X13").Value = -4
$ws.Range("AB13").Value = -14
$ws.Range("X14").Value = -4
$ws.Range("AB14").Value = -15
$ws.Range("X15").Value = -4
$ws.Range("AB15").Value = -16
$ws.Range("X16").Value = -4
$ws.Range("AB16").Value = -16
$ws.Range("X17").Value = -4
$ws.Range("AB17").Value = -17
$ws.Range("X18").Value = -4
$ws.Range("AB18").Value = -17
$ws.Range("X19").Value = -4
$ws.Range("AB19").Value = -18
$ws.Range("X20").Value = -4
$ws.Range("AB20").Value = -18
$ws.Range("X21").Value = -4
$ws.Range("AB21").Value = -16
$ws.Range("X22").Value = -4
$ws.Range("AB22").Value = -17
$ws.Range("X23").Value = -4
$ws.Range("AB23").Value = -18
$ws.Range("X24").Value = -4
$ws.Range("AB24").Value = -18
$ws.Range("X25").Value = -4
$ws.Range("AB25").Value = -17
$ws.Range("X26").Value = -5
$ws.Range("AB26").Value = -17
$ws.Range("X27").Value = -5
$ws.Range("AB27").Value = -17
$ws.Range("X28").Value = -5
$ws.Range("AB28").Value = -18
$ws.Range("X29").Value = -5
$ws.Range("AB29").Value = -18
$ws.Range("X30").Value = -5
$ws.Range("AB30").Value = -18
$ws.Range("X31").Value = -4
$ws.Range("AB31").Value = -17
$ws.Range("X32").Value = -4
$ws.Range("AB32").Value = -18
$ws.Range("X33").Value = -4
$ws.Range("AB33").Value = -17
$ws.Range("X34").Value = -4
$ws.Range("AB34").Value = -19
$ws.Range("X35").Value = -4
$ws.Range("AB35").Value = -17
$ws.Range("X36").Value = -4
$ws.Range("AB36").Value = -16
$ws.Range("X37").Value = -4
$ws.Range("AB37").Value = -16
$ws.Range("X38").Value = -6
$ws.Range("AB38").Value = -17
$ws.Range("X39").Value = -6
$ws.Range("AB39").Value = -17
$ws.Range("X40").Value = -6
$ws.Range("AB40").Value = -17

Write-Output "Applied all cell updates"
